# HEP_Weekly_6December2019.pptx - "Minor change in x-axis because the
# numbers were not shown finely" + the deck's date fields rolled from
# 12/3/19 to 12/5/19 (the presentation was re-saved two days later, so
# every "datetime1" date placeholder picked up the new cached value).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Every "Date Placeholder" showing the old cached date text gets the
#    new one. (PowerPoint has no automation entry point that lets you
#    rewrite a field's cached text while leaving the <a:fld> wrapper in
#    place - Find/Replace and InsertDateTime both leave field runs
#    untouched - so we just retype the placeholder text directly, which
#    is exactly what happens in the UI when a user overtypes a field.)
# ---------------------------------------------------------------------
$oldDate = "12/3/19"
$newDate = "12/5/19"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                $tr = $shape.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Resize the red highlight rectangle ("Rectangle 42") drawn over the
#    x-axis on slide 5 so its bottom edge lines up with the finer axis
#    labels - only the vertical extent moves, the box gets shorter.
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shape = $slide5.Shapes.Item($i)
    if ($shape.Name -eq "Rectangle 42") {
        $shape.Top = 137.4894
        $shape.Height = 99.13575
    }
}
